# Update the "取得日時" (acquired datetime) timestamps in column A
# for the data rows of the "ランサーズ" sheet from 12:43:52 to 12:58:37
# (same date, 2026-01-22), reflecting a new append run at 12:58 JST.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-22 12:58:37"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
